# Update weekly ranking [2025-10-22]
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap B7 and B8
$b7 = $ws.Range("B7").Value()
$b8 = $ws.Range("B8").Value()
$ws.Range("B7").Value = $b8
$ws.Range("B8").Value = $b7

# Swap B27 and B28
$b27 = $ws.Range("B27").Value()
$b28 = $ws.Range("B28").Value()
$ws.Range("B27").Value = $b28
$ws.Range("B28").Value = $b27

# Rotate B40, B41, B42 up by one (B40<-B41, B41<-B42, B42<-B40)
$b40 = $ws.Range("B40").Value()
$b41 = $ws.Range("B41").Value()
$b42 = $ws.Range("B42").Value()
$ws.Range("B40").Value = $b41
$ws.Range("B41").Value = $b42
$ws.Range("B42").Value = $b40
